$d = $word.ActiveDocument

# The last paragraph in the document is the "HISTORY RESERVATION ..." bullet.
# Insert a brand-new bullet paragraph right after it, inheriting the same
# list/paragraph formatting (Paragrafoelenco style, numId 3), then fill in
# its text.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()

$d.Paragraphs.Last.Range.Text = "NELLA CITTA’ ABBIAMO DETTO USER DEFINIED MA NON COME OTTENIAMO LA CITTA’ (VEDI PEZZO WS DI US"
